# Updating to Version 1.3
# Applies the content edits to the RasterCompare sheet of the
# FFRMS_RasterQC_Configuration workbook:
#  - A6: fix typo FFRNS -> FFRMS
#  - B6: clear the 0.2% raster example path (now left blank by default)
#  - C6: add guidance note about leaving the 0.2% raster blank
#  - B7: rename the sample output spreadsheet name to a generic XXX_ prefix
#  - C2: extend folder-path guidance note to mention "for rasters"
#  - Row 2: remove the custom (wrapped) row height, back to default
#  - Selection: move the active selection to B8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "FFRMS 0.2% ACF raster"
$ws.Range("B6").ClearContents() | Out-Null
$ws.Range("C6").Value = "Leave it blank if there is no 0.2% raster in this project"
$ws.Range("B7").Value = "XXX_Riverine_Raster_QC_Result"
$ws.Range("C2").Value = "Please include full folder path, raster name and extention(.tif) for rasters"

# Row 2 no longer needs the taller wrapped height -> restore auto height.
$ws.Rows(2).AutoFit() | Out-Null

# Move the saved selection to B8 (matches the saved workbook view state).
$ws.Range("B8").Select() | Out-Null
